$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Respuesta" column (B). The old "Estado" column (C)
# shifts left to become the new column B.
$ws.Columns.Item(2).Delete()

# The comparison no longer cares about answer order, so the previously
# mismatched row now evaluates as correct.
$ws.Range("B4").Value = "Correcto"
